# ----------------------------------------------------------------------------
# Scheduled market-data refresh: recompute Leve profitability numbers
#
# The workbook tracks, per Leve (per job sheet), the current Eorzea market
# prices pulled for each turn-in item (currentAveragePrice / NQ / HQ) and the
# resulting total price / profit for turning the Leve in (LevePriceNQ/HQ,
# LeveProfitNQ/HQ). This run refreshes those price-derived columns (H:N) with
# newly-sampled market-board data for the rows whose source item price moved
# since the last refresh. Rows/columns not listed are unaffected.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1548.5818
$ws.Range("I15").Value = 1548.5818
$ws.Range("K15").Value = 4645.7454
$ws.Range("M15").Value = -4476.7454
$ws.Range("H55").Value = 487.85715
$ws.Range("I55").Value = 479
$ws.Range("J55").Value = 499.66666
$ws.Range("K55").Value = 479
$ws.Range("L55").Value = 499.66666
$ws.Range("M55").Value = -265
$ws.Range("N55").Value = -927.66666
$ws.Range("H76").Value = 4139.2
$ws.Range("I76").Value = 2900
$ws.Range("J76").Value = 4449
$ws.Range("K76").Value = 2900
$ws.Range("L76").Value = 4449
$ws.Range("M76").Value = -2585
$ws.Range("N76").Value = -5079
$ws.Range("H79").Value = 4139.2
$ws.Range("I79").Value = 2900
$ws.Range("J79").Value = 4449
$ws.Range("K79").Value = 2900
$ws.Range("L79").Value = 4449
$ws.Range("M79").Value = -1808
$ws.Range("N79").Value = -6633
$ws.Range("H80").Value = 14707496
$ws.Range("I80").Value = 19232576
$ws.Range("J80").Value = 989.5
$ws.Range("K80").Value = 57697728
$ws.Range("L80").Value = 2968.5
$ws.Range("M80").Value = -57696730
$ws.Range("N80").Value = -4964.5
$ws.Range("H83").Value = 14707496
$ws.Range("I83").Value = 19232576
$ws.Range("J83").Value = 989.5
$ws.Range("K83").Value = 173093184
$ws.Range("L83").Value = 8905.5
$ws.Range("M83").Value = -173088192
$ws.Range("N83").Value = -18889.5
$ws.Range("H97").Value = 2441.75
$ws.Range("J97").Value = 2441.75
$ws.Range("L97").Value = 7325.25
$ws.Range("N97").Value = -8317.25
$ws.Range("H103").Value = 723.4091
$ws.Range("J103").Value = 675.1875
$ws.Range("L103").Value = 2025.5625
$ws.Range("N103").Value = -3197.5625
$ws.Range("H112").Value = 3994.3928
$ws.Range("J112").Value = 4172.8076
$ws.Range("L112").Value = 12518.4228
$ws.Range("N112").Value = -14734.4228
$ws.Range("H113").Value = 44143.766
$ws.Range("J113").Value = 73194.8
$ws.Range("L113").Value = 73194.8
$ws.Range("N113").Value = -79702.8
$ws.Range("H121").Value = 4632.385
$ws.Range("J121").Value = 4632.385
$ws.Range("L121").Value = 13897.155
$ws.Range("N121").Value = -17391.155
$ws.Range("H125").Value = 1365.9354
$ws.Range("I125").Value = 943.75
$ws.Range("K125").Value = 8493.75
$ws.Range("M125").Value = -6033.75
$ws.Range("H132").Value = 10743.956
$ws.Range("I132").Value = 2895.611
$ws.Range("K132").Value = 8686.832999999999
$ws.Range("M132").Value = -6156.832999999999
$ws.Range("H135").Value = 4569.375
$ws.Range("I135").Value = 1479.4546
$ws.Range("K135").Value = 13315.0914
$ws.Range("M135").Value = -10780.0914
$ws.Range("H137").Value = 2239.5
$ws.Range("I137").Value = 1886.5
$ws.Range("J137").Value = 3298.5
$ws.Range("K137").Value = 5659.5
$ws.Range("L137").Value = 9895.5
$ws.Range("M137").Value = -3109.5
$ws.Range("N137").Value = -14995.5
$ws.Range("H140").Value = 64574.555
$ws.Range("J140").Value = 64213.707
$ws.Range("L140").Value = 64213.707
$ws.Range("N140").Value = -74573.70699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2646.45
$ws.Range("I45").Value = 2006.1538
$ws.Range("K45").Value = 2006.1538
$ws.Range("M45").Value = -1629.1538
$ws.Range("H61").Value = 5790.7144
$ws.Range("J61").Value = 7880.125
$ws.Range("L61").Value = 7880.125
$ws.Range("N61").Value = -8304.125
$ws.Range("H74").Value = 1424.0588
$ws.Range("I74").Value = 1183.0834
$ws.Range("K74").Value = 1183.0834
$ws.Range("M74").Value = -309.0834
$ws.Range("H77").Value = 1424.0588
$ws.Range("I77").Value = 1183.0834
$ws.Range("K77").Value = 5915.416999999999
$ws.Range("M77").Value = -1547.416999999999
$ws.Range("H88").Value = 53165.4
$ws.Range("J88").Value = 58961.555
$ws.Range("L88").Value = 58961.555
$ws.Range("N88").Value = -59773.555
$ws.Range("H91").Value = 53165.4
$ws.Range("J91").Value = 58961.555
$ws.Range("L91").Value = 58961.555
$ws.Range("N91").Value = -61769.555
$ws.Range("H110").Value = 1474.5
$ws.Range("I110").Value = 1399.4286
$ws.Range("K110").Value = 1399.4286
$ws.Range("M110").Value = 645.5714
$ws.Range("H132").Value = 31374.416
$ws.Range("I132").Value = 59497.4
$ws.Range("K132").Value = 178492.2
$ws.Range("M132").Value = -175962.2
$ws.Range("H136").Value = 5790.7144
$ws.Range("J136").Value = 7880.125
$ws.Range("L136").Value = 23640.375
$ws.Range("N136").Value = -28740.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1332.5
$ws.Range("J22").Value = 1332
$ws.Range("L22").Value = 1332
$ws.Range("N22").Value = -1678
$ws.Range("H86").Value = 2310.7036
$ws.Range("I86").Value = 1595.0952
$ws.Range("J86").Value = 4815.3335
$ws.Range("K86").Value = 1595.0952
$ws.Range("L86").Value = 4815.3335
$ws.Range("M86").Value = -472.0952
$ws.Range("N86").Value = -7061.3335
$ws.Range("H89").Value = 2310.7036
$ws.Range("I89").Value = 1595.0952
$ws.Range("J89").Value = 4815.3335
$ws.Range("K89").Value = 7975.476
$ws.Range("L89").Value = 24076.6675
$ws.Range("M89").Value = -2359.476
$ws.Range("N89").Value = -35308.6675
$ws.Range("H99").Value = 2208.75
$ws.Range("I99").Value = 1747
$ws.Range("J99").Value = 3132.25
$ws.Range("K99").Value = 1747
$ws.Range("L99").Value = 3132.25
$ws.Range("M99").Value = -249
$ws.Range("N99").Value = -6128.25
$ws.Range("H134").Value = 2717.3333
$ws.Range("I134").Value = 2198.0454
$ws.Range("K134").Value = 6594.1362
$ws.Range("M134").Value = -4059.1362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2466.5
$ws.Range("J58").Value = 3764
$ws.Range("L58").Value = 3764
$ws.Range("N58").Value = -4170
$ws.Range("H62").Value = 40717.637
$ws.Range("I62").Value = 5966
$ws.Range("K62").Value = 5966
$ws.Range("M62").Value = -5342
$ws.Range("H65").Value = 40717.637
$ws.Range("I65").Value = 5966
$ws.Range("K65").Value = 29830
$ws.Range("M65").Value = -26710
$ws.Range("H93").Value = 166681810
$ws.Range("I93").Value = 8975
$ws.Range("J93").Value = 500027500
$ws.Range("K93").Value = 8975
$ws.Range("L93").Value = 500027500
$ws.Range("M93").Value = -7103
$ws.Range("N93").Value = -500031244
$ws.Range("H107").Value = 892.5333000000001
$ws.Range("I107").Value = 531
$ws.Range("K107").Value = 531
$ws.Range("M107").Value = 1389
$ws.Range("H134").Value = 2651.423
$ws.Range("I134").Value = 2500.647
$ws.Range("J134").Value = 2936.2222
$ws.Range("K134").Value = 7501.941
$ws.Range("L134").Value = 8808.6666
$ws.Range("M134").Value = -4966.941
$ws.Range("N134").Value = -13878.6666
$ws.Range("H136").Value = 2466.5
$ws.Range("J136").Value = 3764
$ws.Range("L136").Value = 11292
$ws.Range("N136").Value = -16392

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1448.2858
$ws.Range("I34").Value = 310.72726
$ws.Range("K34").Value = 932.18178
$ws.Range("M34").Value = -848.18178
$ws.Range("H52").Value = 5388
$ws.Range("J52").Value = 5388
$ws.Range("L52").Value = 16164
$ws.Range("N52").Value = -16696

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").ClearContents()
$ws.Range("H94").Value = 40149.332
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H95").Value = 26329.8
$ws.Range("J95").Value = 26329.8
$ws.Range("L95").Value = 26329.8
$ws.Range("N95").Value = -31821.8
$ws.Range("H122").Value = 396923.8
$ws.Range("I122").Value = 734791.6
$ws.Range("K122").Value = 2204374.8
$ws.Range("M122").Value = -2201924.8
$ws.Range("H132").Value = 4239.25
$ws.Range("I132").Value = 3414.2856
$ws.Range("K132").Value = 10242.8568
$ws.Range("M132").Value = -7712.856800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3787.4707
$ws.Range("I132").Value = 2622.7693
$ws.Range("K132").Value = 7868.3079
$ws.Range("M132").Value = -5338.3079
$ws.Range("H136").Value = 3593.2954
$ws.Range("I136").Value = 2458.3225
$ws.Range("J136").Value = 6299.769
$ws.Range("K136").Value = 7374.967500000001
$ws.Range("L136").Value = 18899.307
$ws.Range("M136").Value = -4824.967500000001
$ws.Range("N136").Value = -23999.307

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10536.546
$ws.Range("I81").Value = 4641
$ws.Range("K81").Value = 9282
$ws.Range("M81").Value = -8221
$ws.Range("H84").Value = 10536.546
$ws.Range("I84").Value = 4641
$ws.Range("K84").Value = 46410
$ws.Range("M84").Value = -41106
$ws.Range("H132").Value = 1999
$ws.Range("I132").Value = 964.1539
$ws.Range("J132").Value = 3222
$ws.Range("K132").Value = 2892.4617
$ws.Range("L132").Value = 9666
$ws.Range("M132").Value = -362.4616999999998
$ws.Range("N132").Value = -14726
